# Generate Report for Handoff
# Rewrites the localization-status report rows to reflect a fresh handoff run:
#  - the source .md file was re-uploaded under a new guid
#  - two new dependent .png files were picked up and reported as well

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# New identifiers produced by this handoff run
$mdGuid   = "5b164d0d-a7fd-4161-9ed8-d1b01f1c829b"
$png1Guid = "806acba8-3d75-47d5-a996-8d9878015983"
$png2Guid = "98d7c109-3a7a-46f0-a119-894c466f29b9"

$mdFile   = "$mdGuid.md"
$png1File = "$png1Guid.png"
$png2File = "$png2Guid.png"

$zhXlf = "$mdGuid.b2a0b45c658c1580c639612b26cfd5b428bbb2b2.zh-cn.xlf"
$deXlf = "$mdGuid.b2a0b45c658c1580c639612b26cfd5b428bbb2b2.de-de.xlf"

$png1Target = "1d4f14c9a5afd187cec193b01801856e103842f7.png"
$png2Target = "10cba0734b6389232678fbf36981d240147cf95a.png"

$overviewDate = "2016-48-19 20:48:42"
$zhDate       = "2016-03-19 20:48:39"
$deDate       = "2016-03-19 20:48:42"

$ready   = "Ready for handoff"
$epoch   = "0001-01-01 00:00:00"
$include = "Include"
$isDep   = "IsDependency"
$depFrom = "e2e\$mdFile"

$srcBase    = "https://github.com/OpenLocalizationTest/oltest/blob/3bb437effbc6cc8de3bfa220aa80c05fafbc74b5/e2e"
$zhHtBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd147698c13a893554668d3f2121c1baa0ead7e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHtBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28ee0c7fafe3109a8a7f47cedbe7bd9374bdce9e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$hlColor = 15570276   # BGR-encoded 0x6495ED so the saved RGB is FF6495ED, matching the workbook's HyperLink font

function Format-Link($ws, $cellRef, $text, $url) {
    $ws.Range($cellRef).Value2 = $text
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text) | Out-Null
    $ws.Range($cellRef).Font.Underline = 1
    $ws.Range($cellRef).Font.Color = $hlColor
    $ws.Range($cellRef).Font.Name = "Calibri"
    $ws.Range($cellRef).Font.Size = 11
}

function Set-DateText($ws, $cellRef, $text) {
    $ws.Range($cellRef).Value2 = $text
    $ws.Range($cellRef).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Overview sheet: drop all existing hyperlinks first so stale ones do not
# linger, then rebuild the three data rows (existing row 2 plus two new rows
# for the dependent png files).
# ---------------------------------------------------------------------------
$ws1.Hyperlinks.Delete()

Format-Link $ws1 "A2" $mdFile "$srcBase/$mdFile"
$ws1.Range("B2").Value2 = $ready
$ws1.Range("C2").Value2 = $ready
$ws1.Range("D2").Value2 = $overviewDate

Format-Link $ws1 "A3" $png1File "$srcBase/$png1File"
$ws1.Range("B3").Value2 = $ready
$ws1.Range("C3").Value2 = $ready
$ws1.Range("D3").Value2 = $overviewDate

Format-Link $ws1 "A4" $png2File "$srcBase/$png2File"
$ws1.Range("B4").Value2 = $ready
$ws1.Range("C4").Value2 = $ready
$ws1.Range("D4").Value2 = $overviewDate

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Hyperlinks.Delete()

Format-Link $ws2 "A2" $mdFile "$srcBase/$mdFile"
Format-Link $ws2 "B2" ".md" "$srcBase/$mdFile"
$ws2.Range("C2").Value2 = $ready
Format-Link $ws2 "D2" $zhXlf "$zhHtBase/$zhXlf"
Set-DateText $ws2 "E2" $zhDate
$ws2.Range("H2").Value2 = $epoch
$ws2.Range("I2").Value2 = $include

Format-Link $ws2 "A3" $png1File "$srcBase/$png1File"
Format-Link $ws2 "B3" ".png" "$srcBase/$png1File"
$ws2.Range("C3").Value2 = $ready
Format-Link $ws2 "D3" $png1Target "$zhHtBase/$png1Target"
Set-DateText $ws2 "E3" $zhDate
$ws2.Range("H3").Value2 = $epoch
$ws2.Range("I3").Value2 = $isDep
$ws2.Range("J3").Value2 = $depFrom

Format-Link $ws2 "A4" $png2File "$srcBase/$png2File"
Format-Link $ws2 "B4" ".png" "$srcBase/$png2File"
$ws2.Range("C4").Value2 = $ready
Format-Link $ws2 "D4" $png2Target "$zhHtBase/$png2Target"
Set-DateText $ws2 "E4" $zhDate
$ws2.Range("H4").Value2 = $epoch
$ws2.Range("I4").Value2 = $isDep
$ws2.Range("J4").Value2 = $depFrom

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Hyperlinks.Delete()

Format-Link $ws3 "A2" $mdFile "$srcBase/$mdFile"
Format-Link $ws3 "B2" ".md" "$srcBase/$mdFile"
$ws3.Range("C2").Value2 = $ready
Format-Link $ws3 "D2" $deXlf "$deHtBase/$deXlf"
Set-DateText $ws3 "E2" $deDate
$ws3.Range("H2").Value2 = $epoch
$ws3.Range("I2").Value2 = $include

Format-Link $ws3 "A3" $png1File "$srcBase/$png1File"
Format-Link $ws3 "B3" ".png" "$srcBase/$png1File"
$ws3.Range("C3").Value2 = $ready
Format-Link $ws3 "D3" $png1Target "$deHtBase/$png1Target"
Set-DateText $ws3 "E3" $deDate
$ws3.Range("H3").Value2 = $epoch
$ws3.Range("I3").Value2 = $isDep
$ws3.Range("J3").Value2 = $depFrom

Format-Link $ws3 "A4" $png2File "$srcBase/$png2File"
Format-Link $ws3 "B4" ".png" "$srcBase/$png2File"
$ws3.Range("C4").Value2 = $ready
Format-Link $ws3 "D4" $png2Target "$deHtBase/$png2Target"
Set-DateText $ws3 "E4" $deDate
$ws3.Range("H4").Value2 = $epoch
$ws3.Range("I4").Value2 = $isDep
$ws3.Range("J4").Value2 = $depFrom

Write-Host "Overview hyperlinks:" $ws1.Hyperlinks.Count
Write-Host "zh-cn hyperlinks:" $ws2.Hyperlinks.Count
Write-Host "de-de hyperlinks:" $ws3.Hyperlinks.Count
